$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) cell format, used to restore
# the original style on cells that we briefly mark as Text so Excel keeps
# numeric-looking strings (e.g. "60.00") verbatim instead of normalizing them.
$plainStyle = $ws.Range('B2').Style

$ws.Range('D2').Value = '37.779.62'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.042.60'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.51'
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.00'
$ws.Range('D7').Style = $plainStyle
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0843'
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = '  +3.13%  '
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '2.346.01'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.35'
$ws.Range('D13').Style = $plainStyle
$ws.Range('E13').Value = '  -2.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.10'
$ws.Range('D14').Style = $plainStyle
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.47'
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = '  +5.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.763'
$ws.Range('D16').Style = $plainStyle
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '2.036.13'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '37.733.15'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.36'
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.93'
$ws.Range('D20').Style = $plainStyle
$ws.Range('E20').Value = '  -2.46%  '
$ws.Range('D21').Value = '0.0₃0827'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  +3.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.94'
$ws.Range('D26').Style = $plainStyle
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.35'
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('E32').Value = '  +8.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.36'
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.49'
$ws.Range('D34').Style = $plainStyle
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.61'
$ws.Range('D36').Style = $plainStyle
$ws.Range('E36').Value = '  +3.30%  '
$ws.Range('E37').Value = '  +3.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.45'
$ws.Range('D38').Style = $plainStyle
$ws.Range('E38').Value = '  +6.03%  '
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.92'
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = '  +5.50%  '
$ws.Range('D41').Value = '1.535.29'
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.75'
$ws.Range('D42').Style = $plainStyle
$ws.Range('E42').Value = '  +0.53%  '
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0903'
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.14'
$ws.Range('D46').Style = $plainStyle
$ws.Range('E46').Value = '  +5.53%  '
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.01'
$ws.Range('D50').Style = $plainStyle
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('D51').Value = '2.234.22'
$ws.Range('E51').Value = '  +0.24%  '
